# Bugfixed evaluation and simulated rt_data for components
# Update the y_1 (column C) and y_1_forecast (column E) naive-forecaster values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 6.109216616889168
$ws.Range("E2").Value = 8.045645122021906
$ws.Range("C3").Value = -5.232639093663815
$ws.Range("E3").Value = -0.9943400519801915
$ws.Range("C4").Value = 2.12454018480297
$ws.Range("E4").Value = 1.203634802640963
$ws.Range("C5").Value = 8.081020954067775
$ws.Range("E5").Value = 2.257871268432821
$ws.Range("C6").Value = 4.489210662380949
$ws.Range("E6").Value = 7.472658273721078
$ws.Range("C7").Value = -0.8752093743685352
$ws.Range("E7").Value = 1.998870338019265
$ws.Range("C8").Value = 4.891728508251214
$ws.Range("E8").Value = 2.926340920335191
$ws.Range("C9").Value = 4.818339085077583
$ws.Range("E9").Value = 4.241902819910548
$ws.Range("C10").Value = 4.067959312311897
$ws.Range("E10").Value = 4.36243732366437
$ws.Range("C11").Value = 4.613634856640747
$ws.Range("E11").Value = 4.246555741688218
$ws.Range("C12").Value = 3.898744563937395
$ws.Range("E12").Value = 3.529300656691237
$ws.Range("C13").Value = 3.841510956591465
$ws.Range("E13").Value = 4.83848589746565
$ws.Range("C14").Value = -0.8225206269755425
$ws.Range("E14").Value = 1.159653508089242
$ws.Range("C15").Value = 4.409066926520455
$ws.Range("E15").Value = 3.034999751677669
$ws.Range("C16").Value = 8.064077385547574
$ws.Range("E16").Value = 3.690055931494096
$ws.Range("C17").Value = 0.3551698673347259
$ws.Range("E17").Value = 4.038118345571751
$ws.Range("C18").Value = -2.267078452724969
$ws.Range("E18").Value = 0.5378929214800987
$ws.Range("C19").Value = 1.326993065386817
$ws.Range("E19").Value = -0.2087957186147071
